$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.939.39"
$ws.Range("E2").Value = "'  +7.79%  "
$ws.Range("D3").Value = "'1.818.56"
$ws.Range("E3").Value = "'  +5.12%  "
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("D5").Value = "'245.82"
$ws.Range("E5").Value = "'  +2.20%  "
$ws.Range("D7").Value = "'0.4926"
$ws.Range("E7").Value = "'  +1.83%  "
$ws.Range("D8").Value = "'44.37"
$ws.Range("E8").Value = "'  +7.14%  "
$ws.Range("D9").Value = "'0.2765"
$ws.Range("E9").Value = "'  +6.33%  "
$ws.Range("D10").Value = "'0.06376"
$ws.Range("E10").Value = "'  +3.19%  "
$ws.Range("D11").Value = "'1.819.07"
$ws.Range("E11").Value = "'  +5.16%  "
$ws.Range("D12").Value = "'16.64"
$ws.Range("E12").Value = "'  +4.24%  "
$ws.Range("D13").Value = "'0.07045"
$ws.Range("E13").Value = "'  +2.60%  "
$ws.Range("B14").Value = "'Polygon"
$ws.Range("C14").Value = "'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.6424"
$ws.Range("E14").Value = "'  +6.50%  "
$ws.Range("B15").Value = "'Litecoin"
$ws.Range("C15").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'83.95"
$ws.Range("E15").Value = "'  +9.06%  "
$ws.Range("D16").Value = "'4.684"
$ws.Range("E16").Value = "'  +4.99%  "
$ws.Range("D17").Value = "'28.955.03"
$ws.Range("E17").Value = "'  +7.92%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("E18").Value = "'  +0.12%  "
$ws.Range("B19").Value = "'BinanceUSD"
$ws.Range("C19").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "'  +0.13%  "
$ws.Range("B20").Value = "'ShibaInu"
$ws.Range("C20").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.000007283"
$ws.Range("E20").Value = "'  +2.05%  "
$ws.Range("D21").Value = "'12.18"
$ws.Range("E21").Value = "'  +7.24%  "
$ws.Range("D22").Value = "'2.052.79"
$ws.Range("E22").Value = "'  +5.17%  "
$ws.Range("D23").Value = "'4.543"
$ws.Range("E23").Value = "'  +3.22%  "
$ws.Range("E24").Value = "'  +4.55%  "
$ws.Range("E25").Value = "'  +5.76%  "
$ws.Range("D26").Value = "'143.89"
$ws.Range("E26").Value = "'  +2.66%  "
$ws.Range("D27").Value = "'128.17"
$ws.Range("E27").Value = "'  +20.35%  "
$ws.Range("D28").Value = "'16.33"
$ws.Range("E28").Value = "'  +7.46%  "
$ws.Range("D29").Value = "'1.877"
$ws.Range("D30").Value = "'1.395"
$ws.Range("E30").Value = "'  +1.43%  "
$ws.Range("D31").Value = "'4.119"
$ws.Range("E31").Value = "'  +4.16%  "
$ws.Range("D32").Value = "'0.08339"
$ws.Range("E32").Value = "'  +5.35%  "
$ws.Range("D33").Value = "'3.764"
$ws.Range("E33").Value = "'  +2.83%  "
$ws.Range("D34").Value = "'0.04943"
$ws.Range("E34").Value = "'  +8.43%  "
$ws.Range("D35").Value = "'1.095"
$ws.Range("E35").Value = "'  +9.47%  "
$ws.Range("D36").Value = "'2.701"
$ws.Range("E36").Value = "'  +4.05%  "
$ws.Range("D37").Value = "'0.6683"
$ws.Range("E37").Value = "'  +8.37%  "
$ws.Range("D38").Value = "'2.276"
$ws.Range("E38").Value = "'  +14.27%  "
$ws.Range("D39").Value = "'2.667"
$ws.Range("E39").Value = "'  +8.77%  "
$ws.Range("D40").Value = "'0.9437"
$ws.Range("E40").Value = "'  +2.35%  "
$ws.Range("D41").Value = "'6.172"
$ws.Range("E41").Value = "'  +9.17%  "
$ws.Range("D42").Value = "'0.01585"
$ws.Range("E42").Value = "'  +6.06%  "
$ws.Range("E43").Value = "'  +0.19%  "
$ws.Range("D44").Value = "'100.89"
$ws.Range("E44").Value = "'  +0.99%  "
$ws.Range("D45").Value = "'0.4055"
$ws.Range("E45").Value = "'  +5.91%  "
$ws.Range("D46").Value = "'7.148"
$ws.Range("E46").Value = "'  +5.31%  "
$ws.Range("D47").Value = "'0.1216"
$ws.Range("E47").Value = "'  +5.16%  "
$ws.Range("D48").Value = "'0.05524"
$ws.Range("E48").Value = "'  +3.01%  "
$ws.Range("B49").Value = "'Elrond"
$ws.Range("C49").Value = "'https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'31.68"
$ws.Range("E49").Value = "'  +5.21%  "
$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.105"
$ws.Range("E50").Value = "'  +2.32%  "
$ws.Range("D51").Value = "'1.297"
$ws.Range("E51").Value = "'  +4.50%  "
